# Update Work Week and Social Spending
# Refreshes the "GDP per Capita" series for Bahrain (rows 2-62, years 1950-2010)
# with revised figures, and appends six new years (2011-2016) of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Revised GDP per Capita values for years 1950-2010 (rows 2-62).
$values = @(
    "3354",
    "3483",
    "3614",
    "3747",
    "3883",
    "4014",
    "4141",
    "4262",
    "4366",
    "4457",
    "4532",
    "4594",
    "4672",
    "4772",
    "4898",
    "5058",
    "5233",
    "5423",
    "5611",
    "5812",
    "6038",
    "6349",
    "6692",
    "6975",
    "7291",
    "6252",
    "6875",
    "7084",
    "7037",
    "6730",
    "6994",
    "6875",
    "7036",
    "7240",
    "7264",
    "6972",
    "6880",
    "6786",
    "6802",
    "6751",
    "6542",
    "7059.95687227835",
    "7993.89275633257",
    "9102.06612385235",
    "9935.00393029083",
    "10695.3924335854",
    "11640.8213531142",
    "12528.3895189374",
    "13776.0082816373",
    "15263.7574913346",
    "17021.7223794001",
    "18152.9626577984",
    "19487.8412319234",
    "21391.6119255483",
    "23389.602406732",
    "25230.9000521677",
    "26840.124270106",
    "28884.6973671033",
    "30612.4978133559",
    "31729.5663186993",
    "34057.7132274499"
)

# New rows for years 2011-2016 (rows 63-68).
$newYears = @(2011, 2012, 2013, 2014, 2015, 2016)
$newValues = @("36372", "36949", "38493", "39799", "40483", "41078")

$firstRow = 2
$lastExistingRow = $firstRow + $values.Count - 1          # 62
$lastNewRow = $lastExistingRow + $newYears.Count           # 68

# The "Data" column stores these numeric-looking figures as text in the
# source workbook. Temporarily format column E as Text across the whole
# range (existing + new rows) so the values are written as text rather than
# being reinterpreted as numbers, then clear the formatting again so the
# cells fall back to the default (unstyled) look, matching the source.
$dataRangeAddress = "E$firstRow`:E$lastNewRow"
$ws.Range($dataRangeAddress).NumberFormat = "@"

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $firstRow + $i
    $ws.Range("E$row").Value = $values[$i]
}

for ($i = 0; $i -lt $newYears.Count; $i++) {
    $row = $lastExistingRow + 1 + $i
    $ws.Range("A$row").Value = 48
    $ws.Range("B$row").Value = "Bahrain"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $newYears[$i]
    $ws.Range("E$row").Value = $newValues[$i]
}

$ws.Range($dataRangeAddress).ClearFormats()
